# Update the "cat table" (categories) in SimpleTest.xlsx:
#  - Column B now holds multi-line "Categories" values (newline separated)
#    instead of the single-line "Content Provenance" values.
#  - Row heights for rows 2, 4 and 7 grow to fit the extra lines of text.
#  - Column widths are set to fit the (now wider/taller) table.
#  - Selection moves to B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$nl = [char]10

# --- Column B values -------------------------------------------------
# Row 1 header keeps its existing value/style (s1, matches A1/D1), just
# re-asserted for clarity.
$ws.Range("B1").Value = "Categories"

$ws.Range("B2").Value = "Content Provenance" + $nl + "Trust and Authenticity" + $nl + "Watermarking"
$ws.Range("B3").Value = "Content Provenance" + $nl + "Trust and Authenticity"
$ws.Range("B4").Value = "Content Provenance" + $nl + "Trust and Authenticity" + $nl + "Asset Identifiers" + $nl + "Rights Declarations" + $nl + "Watermarking"
$ws.Range("B5").Value = "Trust and Authenticity"
$ws.Range("B6").Value = "Watermarking"
$ws.Range("B7").Value = "Content Provenance" + $nl + "Trust and Authenticity" + $nl + "Rights Declarations"

# --- Column B styles ---------------------------------------------------
# B1 -> same style as A1 (s=1)
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# B2, B3, B4, B7 -> same style as A2.. (s=4, wraps text, no border)
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("B7").PasteSpecial(-4122)

# B5, B6 -> same style as C2.. (s=5, no wrap, no border)
$ws.Range("C2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Row heights ---------------------------------------------------
$ws.Rows.Item(2).RowHeight = 80
$ws.Rows.Item(4).RowHeight = 120
$ws.Rows.Item(7).RowHeight = 380

# --- Column widths ---------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.330729166666668
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(3).ColumnWidth = 46.998697916666664
$ws.Columns.Item(4).ColumnWidth = 21.166666666666668
$ws.Columns.Item(5).ColumnWidth = 78.16666666666667
$ws.Columns.Item(6).ColumnWidth = 10.330729166666666
$ws.Columns.Item(7).ColumnWidth = 26.330729166666668
$ws.Columns.Item(8).ColumnWidth = 10.166666666666666

# --- Selection ---------------------------------------------------
$ws.Range("B1").Select()
